$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = 'Background
The outbreak of a novel coronavirus (2019-nCoV)-infected pneumonia (NCIP) is currently ongoing in China.
 Most of the critically ill patients received high-flow nasal cannula (HFNC) oxygen therapy.
 However, the experience of HFNC in this population is lacking.
Methods
We retrospectively screened 318 confirmed patients with NCIP in two hospitals of Chongqing, China, from January 1st to March 4th, 2020. Among them, 27 (8.4%) patients experienced severe acute respiratory failure including 17 patients (63%) treated with HFNC as first-line therapy, 9 patients (33%) treated with noninvasive ventilation (NIV) and one patient (4%) treated with invasive ventilation.
 HFNC failure was defined by the need of NIV or intubation as rescue therapy.
Results
Of the 17 HFNC patients, 7 (41%) experienced HFNC failure.
 The HFNC failure rate was 0% (0/6) in patients with PaO2/FiO2 &gt; 200 mm Hg vs.
 63% (7/11) in those with PaO2/FiO2 ≤ 200 mm Hg (p = 0.04).
 Compared with baseline data, the respiratory rate significantly decreased after 1–2 h of HFNC in successful group [median 26 (IQR: 25–29) vs.
 23 (22–25), p = 0.03].
 However, it did not in the unsuccessful group.
 After initiation of NIV as rescue therapy among the 7 patients with HFNC failure, PaO2/FiO2 significantly improved after 1–2 h of NIV [median 172 (150–208) mmHg vs.
 114 (IQR: 79–130) under HFNC, p = 0.04].
 However, two out of seven (29%) patients with NIV as rescue therapy ultimately received intubation.
 Among the 27 patients with severe acute respiratory failure, four patients were eventually intubated (15%).
Conclusions
Our study indicated that HFNC was the most common ventilation support for patients with NCIP.
 Patients with lower PaO2/FiO2 were more likely to experience HFNC failure.
'
$authors = '[Ke%Wang%NULL%0,      Wei%Zhao%NULL%1,      Ji%Li%NULL%1,      Weiwei%Shu%shuweiwei361@163.com%1,      Jun%Duan%duanjun412589@163.com%1]'

$ws.Range("D16").Value = $abstract
$ws.Range("E16").Value = $authors
